$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2368.0312
$ws.Range("J80").Value = 2980.5217
$ws.Range("L80").Value = 8941.5651
$ws.Range("N80").Value = -10937.5651

$ws.Range("H83").Value = 2368.0312
$ws.Range("J83").Value = 2980.5217
$ws.Range("L83").Value = 26824.6953
$ws.Range("N83").Value = -36808.6953

$ws.Range("H94").Value = 658.44446
$ws.Range("I94").Value = 658.44446
$ws.Range("K94").Value = 658.44446
$ws.Range("M94").Value = -207.44446

$ws.Range("H95").Value = 31156
$ws.Range("J95").Value = 31156
$ws.Range("L95").Value = 31156
$ws.Range("N95").Value = -36648

$ws.Range("H113").Value = 12663.883
$ws.Range("I113").Value = 5978.8335
$ws.Range("K113").Value = 5978.8335
$ws.Range("M113").Value = -2724.8335

$ws.Range("H121").Value = 3180.8
$ws.Range("J121").Value = 3180.8
$ws.Range("L121").Value = 9542.400000000001
$ws.Range("N121").Value = -13036.4

$ws.Range("H132").Value = 16890.121
$ws.Range("I132").Value = 1215.8
$ws.Range("J132").Value = 173633.33
$ws.Range("K132").Value = 3647.4
$ws.Range("L132").Value = 520899.99
$ws.Range("M132").Value = -1117.4
$ws.Range("N132").Value = -525959.99

$ws.Range("H140").Value = 37490.9
$ws.Range("J140").Value = 39888.89
$ws.Range("L140").Value = 39888.89
$ws.Range("N140").Value = -50248.89

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1015
$ws.Range("I21").Value = 1015
$ws.Range("K21").Value = 1015
$ws.Range("M21").Value = -641

$ws.Range("H88").Value = 1777.8182
$ws.Range("I88").Value = 1957.6666
$ws.Range("J88").Value = 1710.375
$ws.Range("K88").Value = 1957.6666
$ws.Range("L88").Value = 1710.375
$ws.Range("M88").Value = -1551.6666
$ws.Range("N88").Value = -2522.375

$ws.Range("H91").Value = 1777.8182
$ws.Range("I91").Value = 1957.6666
$ws.Range("J91").Value = 1710.375
$ws.Range("K91").Value = 1957.6666
$ws.Range("L91").Value = 1710.375
$ws.Range("M91").Value = -553.6666
$ws.Range("N91").Value = -4518.375

$ws.Range("H110").Value = 3459.5293
$ws.Range("I110").Value = 3374.1333
$ws.Range("J110").Value = 4100
$ws.Range("K110").Value = 3374.1333
$ws.Range("L110").Value = 4100
$ws.Range("M110").Value = -1329.1333
$ws.Range("N110").Value = -8190

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 55699444
$ws.Range("I86").Value = 62655624
$ws.Range("J86").Value = 50000
$ws.Range("K86").Value = 62655624
$ws.Range("L86").Value = 50000
$ws.Range("M86").Value = -62654501
$ws.Range("N86").Value = -52246

$ws.Range("H89").Value = 55699444
$ws.Range("I89").Value = 62655624
$ws.Range("J89").Value = 50000
$ws.Range("K89").Value = 313278120
$ws.Range("L89").Value = 250000
$ws.Range("M89").Value = -313272504
$ws.Range("N89").Value = -261232

$ws.Range("H107").Value = 9988
$ws.Range("I107").Value = 10753.2
$ws.Range("J107").Value = 8712.666999999999
$ws.Range("K107").Value = 10753.2
$ws.Range("L107").Value = 8712.666999999999
$ws.Range("M107").Value = -8833.200000000001
$ws.Range("N107").Value = -12552.667

$ws.Range("H134").Value = 440
$ws.Range("I134").Value = 440
$ws.Range("K134").Value = 1320
$ws.Range("M134").Value = 1215

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2603.842
$ws.Range("I31").Value = 2248.2
$ws.Range("J31").Value = 3937.5
$ws.Range("K31").Value = 2248.2
$ws.Range("L31").Value = 3937.5
$ws.Range("M31").Value = -1953.2
$ws.Range("N31").Value = -4527.5

$ws.Range("H34").Value = 2603.842
$ws.Range("I34").Value = 2248.2
$ws.Range("J34").Value = 3937.5
$ws.Range("K34").Value = 2248.2
$ws.Range("L34").Value = 3937.5
$ws.Range("M34").Value = -2046.2
$ws.Range("N34").Value = -4341.5

$ws.Range("H86").Value = 10664.471
$ws.Range("I86").Value = 3440
$ws.Range("J86").Value = 17086.223
$ws.Range("K86").Value = 3440
$ws.Range("L86").Value = 17086.223
$ws.Range("M86").Value = -2317
$ws.Range("N86").Value = -19332.223

$ws.Range("H89").Value = 10664.471
$ws.Range("I89").Value = 3440
$ws.Range("J89").Value = 17086.223
$ws.Range("K89").Value = 17200
$ws.Range("L89").Value = 85431.11500000001
$ws.Range("M89").Value = -11584
$ws.Range("N89").Value = -96663.11500000001

$ws.Range("H122").Value = 4800.476
$ws.Range("I122").Value = 3320.4443
$ws.Range("J122").Value = 5910.5
$ws.Range("K122").Value = 9961.332900000001
$ws.Range("L122").Value = 17731.5
$ws.Range("M122").Value = -7511.332900000001
$ws.Range("N122").Value = -22631.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 513.9355
$ws.Range("I12").Value = 358.2857
$ws.Range("J12").Value = 559.3333
$ws.Range("K12").Value = 1074.8571
$ws.Range("L12").Value = 1677.9999
$ws.Range("M12").Value = -901.8571000000002
$ws.Range("N12").Value = -2023.9999

$ws.Range("H80").Value = 7666.6665
$ws.Range("I80").Value = 7250
$ws.Range("K80").Value = 21750
$ws.Range("M80").Value = -20814

$ws.Range("H83").Value = 7666.6665
$ws.Range("I83").Value = 7250
$ws.Range("K83").Value = 65250
$ws.Range("M83").Value = -60570

$ws.Range("H95").Value = 4990
$ws.Range("I95").Value = 4990
$ws.Range("K95").Value = 14970
$ws.Range("M95").Value = -12911

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 58947.285
$ws.Range("I70").Value = 85092.28999999999
$ws.Range("K70").Value = 85092.28999999999
$ws.Range("M70").Value = -84822.28999999999

$ws.Range("H73").Value = 58947.285
$ws.Range("I73").Value = 85092.28999999999
$ws.Range("K73").Value = 85092.28999999999
$ws.Range("M73").Value = -84156.28999999999

$ws.Range("H80").Value = 50248
$ws.Range("I80").Value = 142375.5
$ws.Range("J80").Value = 4184.25
$ws.Range("K80").Value = 142375.5
$ws.Range("L80").Value = 4184.25
$ws.Range("M80").Value = -141377.5
$ws.Range("N80").Value = -6180.25

$ws.Range("H83").Value = 50248
$ws.Range("I83").Value = 142375.5
$ws.Range("J83").Value = 4184.25
$ws.Range("K83").Value = 711877.5
$ws.Range("L83").Value = 20921.25
$ws.Range("M83").Value = -706885.5
$ws.Range("N83").Value = -30905.25

$ws.Range("H102").Value = 3964
$ws.Range("I102").Value = 3669.3333
$ws.Range("K102").Value = 3669.3333
$ws.Range("M102").Value = -2047.3333

$ws.Range("H120").Value = 49910.75
$ws.Range("J120").Value = 49910.75
$ws.Range("L120").Value = 49910.75
$ws.Range("N120").Value = -59586.75

$ws.Range("H132").Value = 7548.15
$ws.Range("I132").Value = 6843.4414
$ws.Range("K132").Value = 20530.3242
$ws.Range("M132").Value = -18000.3242

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 110000
$ws.Range("I40").Value = 20000
$ws.Range("K40").Value = 20000
$ws.Range("M40").Value = -19864

$ws.Range("H46").Value = 2794.5715
$ws.Range("J46").Value = 3588.2222
$ws.Range("L46").Value = 3588.2222
$ws.Range("N46").Value = -3964.2222

$ws.Range("H122").Value = 4697.5625
$ws.Range("J122").Value = 5533.8184
$ws.Range("L122").Value = 16601.4552
$ws.Range("N122").Value = -21501.4552

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = $null

$ws.Range("H122").Value = 2100
$ws.Range("I122").Value = 2100
$ws.Range("K122").Value = 6300
$ws.Range("M122").Value = -3850

$ws.Range("H132").Value = 1855.6875
$ws.Range("I132").Value = 1706.4828
$ws.Range("K132").Value = 5119.4484
$ws.Range("M132").Value = -2589.4484
